$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.269.18'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '2.999.97'
$ws.Range('E3').Value = '  +1.26%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '355.05'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.99'
$ws.Range('E6').Value = '  -2.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.565'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.624'
$ws.Range('E9').Value = '  -1.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.58'
$ws.Range('E10').Value = '  -3.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0861'
$ws.Range('E12').Value = '  -3.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.38'
$ws.Range('E13').Value = '  -2.91%  '
$ws.Range('D14').Value = '3.467.67'
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.81'
$ws.Range('E15').Value = '  -3.14%  '
$ws.Range('D16').Value = '3.000.73'
$ws.Range('E16').Value = '  +1.11%  '
$ws.Range('E17').Value = '  +2.65%  '
$ws.Range('D18').Value = '52.360.40'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('E19').Value = '  +6.02%  '
$ws.Range('E20').Value = '  -2.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.72'
$ws.Range('E21').Value = '  -5.33%  '
$ws.Range('D22').Value = '0.0₃0975'
$ws.Range('E22').Value = '  -1.54%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.68'
$ws.Range('E23').Value = '  -2.55%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '264.87'
$ws.Range('E24').Value = '  -2.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.75'
$ws.Range('E25').Value = '  -2.02%  '
$ws.Range('E26').Value = '  -1.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.03'
$ws.Range('E27').Value = '  -1.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.67'
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  -1.23%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.44'
$ws.Range('E31').Value = '  +1.36%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.38'
$ws.Range('E32').Value = '  -3.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '36.78'
$ws.Range('E33').Value = '  -2.64%  '
$ws.Range('E34').Value = '  +12.25%  '
$ws.Range('E35').Value = '  -4.20%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.22'
$ws.Range('E38').Value = '  -5.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.99'
$ws.Range('E39').Value = '  -5.42%  '
$ws.Range('E40').Value = '  -4.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.72'
$ws.Range('E41').Value = '  +0.95%  '
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.98'
$ws.Range('E43').Value = '  -3.43%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '124.27'
$ws.Range('E44').Value = '  +8.87%  '
$ws.Range('E45').Value = '  -1.65%  '
$ws.Range('D46').Value = '2.129.10'
$ws.Range('E46').Value = '  -2.33%  '
$ws.Range('E47').Value = '  -4.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.41'
$ws.Range('E48').Value = '  -5.12%  '
$ws.Range('E49').Value = '  +2.14%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0336'
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.933'
$ws.Range('E51').Value = '  -0.70%  '
